$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column D ("lesion_volume") entirely; subsequent columns (WPM, FPM) shift left
$ws.Range("D1:D52").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
